# Mise a jour mapping de la partie corps
# 1) Bump the "Date" metadata value on the "Metadata" sheet.
# 2) Insert a new top-level mapping row ("FRCDAImageIllustrative" ->
#    "FRMediaDocument") at the top of the "Mapping Table 1" sheet's data
#    rows, pushing every existing mapping row down by one.

$wb = $excel.ActiveWorkbook

# --- 1) Update the Date value on the Metadata sheet -----------------------
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B8").Value = "2026-01-26T10:27:23+00:00"

# --- 2) Insert a new row at the top of the mapping table on "Mapping Table 1"
$ws = $wb.Worksheets.Item("Mapping Table 1")

$lastRow = 15
$newLastRow = $lastRow + 1

# Extend the formatting one row further down (copy row 15's format into the
# new row 16) before shuffling values, so every cell keeps style index 2
# instead of picking up a blank/default style.
$ws.Range("A" + $lastRow + ":E" + $lastRow).Copy()
$ws.Range("A" + $newLastRow + ":E" + $newLastRow).PasteSpecial(-4122)

# Shift the existing mapping rows (3..15) down by one row, bottom-up so we
# never clobber a row before it has been read.
for ($r = $lastRow; $r -ge 3; $r--) {
    for ($col = 1; $col -le 5; $col++) {
        $srcCell = $ws.Cells.Item($r, $col)
        $dstCell = $ws.Cells.Item($r + 1, $col)
        $dstCell.Value = $srcCell.Value2
    }
}

# Write the new top-level mapping row into the now-empty row 3.
$ws.Range("A3").Value = "FRCDAImageIllustrative"
$ws.Range("B3").Value = ""
$ws.Range("C3").Value = "equivalent"
$ws.Range("D3").Value = "FRMediaDocument"
$ws.Range("E3").Value = ""

Write-Host "Mapping corps update applied"
